$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S12/G01..G04 task rows (89-96): mark tasks as implemented and
# replace the planning "remarks" text with a short note describing what
# was actually implemented, per the commit message.
$updates = @(
    @{ Row = 89; F = "V1 scoring model implemented (mom20, ATR%, volume) and PRD aligned."; G = "implemented" }
    @{ Row = 90; F = "score_candidate() helper wired into _run_portfolio_simulator with tests."; G = "implemented" }
    @{ Row = 91; F = "Portfolio simulator now supports multiple scored entries per bar under risk constraints."; G = "implemented" }
    @{ Row = 92; F = "Verified cash/equity updates for multi-entry bars; optional caps can be added later."; G = "implemented" }
    @{ Row = 93; F = "Added routing_debug metrics and per-symbol summaries for group backtests."; G = "implemented" }
    @{ Row = 94; F = "Backtests UI shows capital-aware routing debug and explains scoring legend."; G = "implemented" }
    @{ Row = 95; F = "New tests ensure group trades honour risk_config and MIS/CNC rules."; G = "implemented" }
    @{ Row = 96; F = "Regression harness compares single-stock vs group BT behaviour under routing."; G = "implemented" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}
